# Auto-generated script to apply Aegis_Profits data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 308.39285
$ws.Range("I53").Value = 437.25
$ws.Range("J53").Value = 256.85
$ws.Range("K53").Value = 437.25
$ws.Range("L53").Value = 256.85
$ws.Range("M53").Value = 199.75
$ws.Range("N53").Value = -1530.85
$ws.Range("H64").Value = 94070.55
$ws.Range("I64").Value = 335533.34
$ws.Range("J64").Value = 3522
$ws.Range("K64").Value = 335533.34
$ws.Range("L64").Value = 3522
$ws.Range("M64").Value = -335285.34
$ws.Range("N64").Value = -4018
$ws.Range("H67").Value = 94070.55
$ws.Range("I67").Value = 335533.34
$ws.Range("J67").Value = 3522
$ws.Range("K67").Value = 335533.34
$ws.Range("L67").Value = 3522
$ws.Range("M67").Value = -334675.34
$ws.Range("N67").Value = -5238
$ws.Range("H76").Value = 4669.8
$ws.Range("J76").Value = 4837.25
$ws.Range("L76").Value = 4837.25
$ws.Range("N76").Value = -5467.25
$ws.Range("H79").Value = 4669.8
$ws.Range("J79").Value = 4837.25
$ws.Range("L79").Value = 4837.25
$ws.Range("N79").Value = -7021.25
$ws.Range("H99").Value = 11398.556
$ws.Range("I99").Value = 12698.375
$ws.Range("K99").Value = 38095.125
$ws.Range("M99").Value = -36597.125
$ws.Range("H113").Value = 101840
$ws.Range("J113").Value = 1947.5
$ws.Range("L113").Value = 1947.5
$ws.Range("N113").Value = -8455.5
$ws.Range("H129").Value = 8302.214
$ws.Range("J129").Value = 1741.75
$ws.Range("L129").Value = 5225.25
$ws.Range("N129").Value = -15225.25
$ws.Range("H135").Value = 423.5484
$ws.Range("I135").Value = 423.5484
$ws.Range("K135").Value = 3811.9356
$ws.Range("M135").Value = -1276.9356
$ws.Range("H141").Value = 1696.825
$ws.Range("I141").Value = 1231.1818
$ws.Range("J141").Value = 3892
$ws.Range("K141").Value = 3693.5454
$ws.Range("L141").Value = 11676
$ws.Range("M141").Value = 1486.4546
$ws.Range("N141").Value = -22036
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 253220
$ws.Range("I2").Value = 4690
$ws.Range("K2").Value = 4690
$ws.Range("M2").Value = -4577
$ws.Range("H45").Value = 1000000
$ws.Range("I45").Value = 1000000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1000000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -999623
$ws.Range("N45").Value = $null
$ws.Range("H74").Value = 11115057
$ws.Range("I74").Value = 4004
$ws.Range("J74").Value = 16670583
$ws.Range("K74").Value = 4004
$ws.Range("L74").Value = 16670583
$ws.Range("M74").Value = -3130
$ws.Range("N74").Value = -16672331
$ws.Range("H77").Value = 11115057
$ws.Range("I77").Value = 4004
$ws.Range("J77").Value = 16670583
$ws.Range("K77").Value = 20020
$ws.Range("L77").Value = 83352915
$ws.Range("M77").Value = -15652
$ws.Range("N77").Value = -83361651
$ws.Range("H116").Value = 253220
$ws.Range("I116").Value = 4690
$ws.Range("K116").Value = 4690
$ws.Range("M116").Value = -2396
$ws.Range("H122").Value = 1230.2593
$ws.Range("I122").Value = 1143.25
$ws.Range("J122").Value = 1926.3334
$ws.Range("K122").Value = 3429.75
$ws.Range("L122").Value = 5779.0002
$ws.Range("M122").Value = -979.75
$ws.Range("N122").Value = -10679.0002
$ws.Range("H131").Value = 36938
$ws.Range("J131").Value = 36938
$ws.Range("L131").Value = 36938
$ws.Range("N131").Value = -47018
$ws.Range("H134").Value = 65424.5
$ws.Range("J134").Value = 65424.5
$ws.Range("L134").Value = 65424.5
$ws.Range("N134").Value = -75564.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 253220
$ws.Range("I3").Value = 4690
$ws.Range("K3").Value = 4690
$ws.Range("M3").Value = -4576
$ws.Range("H86").Value = 50782.92
$ws.Range("I86").Value = 73081.53
$ws.Range("J86").Value = 3398.375
$ws.Range("K86").Value = 73081.53
$ws.Range("L86").Value = 3398.375
$ws.Range("M86").Value = -71958.53
$ws.Range("N86").Value = -5644.375
$ws.Range("H89").Value = 50782.92
$ws.Range("I89").Value = 73081.53
$ws.Range("J89").Value = 3398.375
$ws.Range("K89").Value = 365407.65
$ws.Range("L89").Value = 16991.875
$ws.Range("M89").Value = -359791.65
$ws.Range("N89").Value = -28223.875
$ws.Range("H99").Value = 1673.8462
$ws.Range("I99").Value = 1542.5
$ws.Range("K99").Value = 1542.5
$ws.Range("M99").Value = -44.5
$ws.Range("H105").Value = 81922.88
$ws.Range("I105").Value = 84915.836
$ws.Range("K105").Value = 84915.836
$ws.Range("M105").Value = -83168.836
$ws.Range("H134").Value = 3474.04
$ws.Range("I134").Value = 3474.04
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10422.12
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7887.119999999999
$ws.Range("N134").Value = $null
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 36803.668
$ws.Range("I31").Value = 1257.125
$ws.Range("J31").Value = 77428.28999999999
$ws.Range("K31").Value = 1257.125
$ws.Range("L31").Value = 77428.28999999999
$ws.Range("M31").Value = -962.125
$ws.Range("N31").Value = -78018.28999999999
$ws.Range("H34").Value = 36803.668
$ws.Range("I34").Value = 1257.125
$ws.Range("J34").Value = 77428.28999999999
$ws.Range("K34").Value = 1257.125
$ws.Range("L34").Value = 77428.28999999999
$ws.Range("M34").Value = -1055.125
$ws.Range("N34").Value = -77832.28999999999
$ws.Range("H94").Value = 1156.6
$ws.Range("I94").Value = 970.6667
$ws.Range("J94").Value = 1203.0834
$ws.Range("K94").Value = 970.6667
$ws.Range("L94").Value = 1203.0834
$ws.Range("M94").Value = -519.6667
$ws.Range("N94").Value = -2105.0834
$ws.Range("H118").Value = 46974.5
$ws.Range("J118").Value = 46974.5
$ws.Range("L118").Value = 46974.5
$ws.Range("N118").Value = -50288.5
$ws.Range("H132").Value = 2804.037
$ws.Range("I132").Value = 2753.1904
$ws.Range("J132").Value = 2982
$ws.Range("K132").Value = 8259.5712
$ws.Range("L132").Value = 8946
$ws.Range("M132").Value = -5729.5712
$ws.Range("N132").Value = -14006
$ws.Range("H135").Value = 43304.668
$ws.Range("I135").Value = 28709
$ws.Range("J135").Value = 45550.152
$ws.Range("K135").Value = 28709
$ws.Range("L135").Value = 45550.152
$ws.Range("M135").Value = -23639
$ws.Range("N135").Value = -55690.152
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 786.89
$ws.Range("I131").Value = 465.08334
$ws.Range("J131").Value = 830.7727
$ws.Range("K131").Value = 1395.25002
$ws.Range("L131").Value = 2492.3181
$ws.Range("M131").Value = 3644.74998
$ws.Range("N131").Value = -12572.3181
$ws.Range("H132").Value = 1836.7273
$ws.Range("J132").Value = 2022.6666
$ws.Range("L132").Value = 18203.9994
$ws.Range("N132").Value = -23263.9994
$ws.Range("H134").Value = 2941.724
$ws.Range("I134").Value = 2935.5557
$ws.Range("K134").Value = 8806.667099999999
$ws.Range("M134").Value = -3736.667099999999
$ws.Range("H137").Value = 60046.668
$ws.Range("I137").Value = 87831.664
$ws.Range("K137").Value = 263494.992
$ws.Range("M137").Value = -258394.992
$ws.Range("H140").Value = 6105.864
$ws.Range("I140").Value = 8316.357
$ws.Range("K140").Value = 24949.071
$ws.Range("M140").Value = -19769.071
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 12598.2
$ws.Range("J48").Value = 12598.2
$ws.Range("L48").Value = 12598.2
$ws.Range("N48").Value = -13568.2
$ws.Range("H80").Value = 250250990
$ws.Range("I80").Value = 500500000
$ws.Range("J80").Value = 1999.5
$ws.Range("K80").Value = 500500000
$ws.Range("L80").Value = 1999.5
$ws.Range("M80").Value = -500499002
$ws.Range("N80").Value = -3995.5
$ws.Range("H83").Value = 250250990
$ws.Range("I83").Value = 500500000
$ws.Range("J83").Value = 1999.5
$ws.Range("K83").Value = 2502500000
$ws.Range("L83").Value = 9997.5
$ws.Range("M83").Value = -2502495008
$ws.Range("N83").Value = -19981.5
$ws.Range("H132").Value = 3533.875
$ws.Range("I132").Value = 3249.4546
$ws.Range("J132").Value = 4159.6
$ws.Range("K132").Value = 9748.363799999999
$ws.Range("L132").Value = 12478.8
$ws.Range("M132").Value = -7218.363799999999
$ws.Range("N132").Value = -17538.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1093.1177
$ws.Range("J22").Value = 1091.0769
$ws.Range("L22").Value = 1091.0769
$ws.Range("N22").Value = -1681.0769
$ws.Range("H27").Value = 1093.1177
$ws.Range("J27").Value = 1091.0769
$ws.Range("L27").Value = 1091.0769
$ws.Range("N27").Value = -1305.0769
$ws.Range("H122").Value = 1722.6
$ws.Range("J122").Value = 2104.5
$ws.Range("L122").Value = 6313.5
$ws.Range("N122").Value = -11213.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 936031.4
$ws.Range("I2").Value = 2015999.4
$ws.Range("J2").Value = 36058
$ws.Range("K2").Value = 2015999.4
$ws.Range("L2").Value = 36058
$ws.Range("M2").Value = -2015887.4
$ws.Range("N2").Value = -36282

Write-Output "Applied $(235) changes"